$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$questions = @(
    @("Q006 ", "What is the name of your favorite childhood hero?"),
    @("Q007", "What is your fathers middle name?"),
    @("Q008", "What is the name of your first school?"),
    @("Q009", "What is the name of your fist crush?"),
    @("Q010", "What is the registration number of your first vehicle?")
)

$row = 6
foreach ($q in $questions) {
    $ws.Cells.Item($row, 1).Value = $q[0]
    $ws.Cells.Item($row, 2).Value = $q[1]
    $row++
}

$ws.Range("A11").Select()
